$d = $word.ActiveDocument

# --- Change 1: "Portfolio:" + " " (two runs) => single run "Portfolio Site: " ---
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t.StartsWith("Portfolio:")) {
        $s = $p.Range.Start
        # First run holds "Portfolio:" (10 chars), the very next run is just a single space.
        $spaceRun = $d.Range($s + 10, $s + 11)
        $spaceRun.Delete()
        $labelRun = $d.Range($s, $s + 10)
        $labelRun.Text = "Portfolio Site: "
        break
    }
}

# --- Change 2: remove the leading "Technical " run from the "Technical Skills" heading ---
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13)
    if ($t -eq "Technical Skills") {
        $s = $p.Range.Start
        $prefixRun = $d.Range($s, $s + 10)
        $prefixRun.Delete()
        break
    }
}
